$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2033.3334
$ws.Range("I62").Value = 987.5
$ws.Range("J62").Value = 4125
$ws.Range("K62").Value = 987.5
$ws.Range("L62").Value = 4125
$ws.Range("M62").Value = -363.5
$ws.Range("N62").Value = -5373

$ws.Range("H65").Value = 2033.3334
$ws.Range("I65").Value = 987.5
$ws.Range("J65").Value = 4125
$ws.Range("K65").Value = 4937.5
$ws.Range("L65").Value = 20625
$ws.Range("M65").Value = -1817.5
$ws.Range("N65").Value = -26865

$ws.Range("H123").Value = 40650
$ws.Range("J123").Value = 40650
$ws.Range("L123").Value = 40650
$ws.Range("N123").Value = -50450

$ws.Range("H137").Value = 4033.7297
$ws.Range("I137").Value = 3473.8572
$ws.Range("J137").Value = 5775.5557
$ws.Range("K137").Value = 10421.5716
$ws.Range("L137").Value = 17326.6671
$ws.Range("M137").Value = -7871.571599999999
$ws.Range("N137").Value = -22426.6671

$ws.Range("H139").Value = 84998.89
$ws.Range("J139").Value = 84998.89
$ws.Range("L139").Value = 84998.89
$ws.Range("N139").Value = -95278.89

$ws.Range("H141").Value = 2059.6924
$ws.Range("I141").Value = 1801.9688
$ws.Range("J141").Value = 3237.8572
$ws.Range("K141").Value = 5405.9064
$ws.Range("L141").Value = 9713.571599999999
$ws.Range("M141").Value = -225.9063999999998
$ws.Range("N141").Value = -20073.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1443
$ws.Range("I45").Value = 1379.2667
$ws.Range("J45").Value = 1682
$ws.Range("K45").Value = 1379.2667
$ws.Range("L45").Value = 1682
$ws.Range("M45").Value = -1002.2667
$ws.Range("N45").Value = -2436

$ws.Range("H74").Value = 1468.4546
$ws.Range("I74").Value = 728.1177
$ws.Range("J74").Value = 3985.6
$ws.Range("K74").Value = 728.1177
$ws.Range("L74").Value = 3985.6
$ws.Range("M74").Value = 145.8823
$ws.Range("N74").Value = -5733.6

$ws.Range("H77").Value = 1468.4546
$ws.Range("I77").Value = 728.1177
$ws.Range("J77").Value = 3985.6
$ws.Range("K77").Value = 3640.5885
$ws.Range("L77").Value = 19928
$ws.Range("M77").Value = 727.4115000000002
$ws.Range("N77").Value = -28664

$ws.Range("H92").Value = 35550
$ws.Range("J92").Value = 35550
$ws.Range("L92").Value = 35550
$ws.Range("N92").Value = -40542

$ws.Range("H122").Value = 5558254.5
$ws.Range("I122").Value = 7694417.5
$ws.Range("J122").Value = 4230
$ws.Range("K122").Value = 23083252.5
$ws.Range("L122").Value = 12690
$ws.Range("M122").Value = -23080802.5
$ws.Range("N122").Value = -17590

$ws.Range("H132").Value = 2016.4166
$ws.Range("I132").Value = 2288.0356
$ws.Range("J132").Value = 1636.15
$ws.Range("K132").Value = 6864.1068
$ws.Range("L132").Value = 4908.450000000001
$ws.Range("M132").Value = -4334.1068
$ws.Range("N132").Value = -9968.450000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1446.0435
$ws.Range("I99").Value = 1366.1875
$ws.Range("J99").Value = 1628.5714
$ws.Range("K99").Value = 1366.1875
$ws.Range("L99").Value = 1628.5714
$ws.Range("M99").Value = 131.8125
$ws.Range("N99").Value = -4624.5714

$ws.Range("H133").Value = 24926.666
$ws.Range("J133").Value = 24926.666
$ws.Range("L133").Value = 24926.666
$ws.Range("N133").Value = -35046.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32218.88
$ws.Range("I31").Value = 43051.76
$ws.Range("J31").Value = 16288.177
$ws.Range("K31").Value = 43051.76
$ws.Range("L31").Value = 16288.177
$ws.Range("M31").Value = -42756.76
$ws.Range("N31").Value = -16878.177

$ws.Range("H34").Value = 32218.88
$ws.Range("I34").Value = 43051.76
$ws.Range("J34").Value = 16288.177
$ws.Range("K34").Value = 43051.76
$ws.Range("L34").Value = 16288.177
$ws.Range("M34").Value = -42849.76
$ws.Range("N34").Value = -16692.177

$ws.Range("H107").Value = 1184.24
$ws.Range("I107").Value = 1572.9375
$ws.Range("J107").Value = 493.22223
$ws.Range("K107").Value = 1572.9375
$ws.Range("L107").Value = 493.22223
$ws.Range("M107").Value = 347.0625
$ws.Range("N107").Value = -4333.22223

$ws.Range("H132").Value = 1325.125
$ws.Range("I132").Value = 1002.2
$ws.Range("J132").Value = 1863.3334
$ws.Range("K132").Value = 3006.6
$ws.Range("L132").Value = 5590.0002
$ws.Range("M132").Value = -476.6000000000004
$ws.Range("N132").Value = -10650.0002

$ws.Range("H134").Value = 2019.7941
$ws.Range("I134").Value = 1612.8214
$ws.Range("J134").Value = 3919
$ws.Range("K134").Value = 4838.4642
$ws.Range("L134").Value = 11757
$ws.Range("M134").Value = -2303.4642
$ws.Range("N134").Value = -16827

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1065.0227
$ws.Range("I5").Value = 250.22223
$ws.Range("K5").Value = 750.66669
$ws.Range("M5").Value = -638.66669

$ws.Range("H33").Value = 700.6667
$ws.Range("I33").Value = 100
$ws.Range("J33").Value = 820.8
$ws.Range("K33").Value = 600
$ws.Range("L33").Value = 4924.799999999999
$ws.Range("M33").Value = -317
$ws.Range("N33").Value = -5490.799999999999

$ws.Range("H44").Value = 942.3
$ws.Range("I44").Value = 303.83334
$ws.Range("J44").Value = 1900
$ws.Range("K44").Value = 911.5000200000001
$ws.Range("L44").Value = 5700
$ws.Range("M44").Value = -513.5000200000001
$ws.Range("N44").Value = -6496

$ws.Range("H64").Value = 3631
$ws.Range("J64").Value = 6500
$ws.Range("L64").Value = 19500
$ws.Range("N64").Value = -20040

$ws.Range("H67").Value = 3631
$ws.Range("J67").Value = 6500
$ws.Range("L67").Value = 19500
$ws.Range("N67").Value = -21372

$ws.Range("H132").Value = 1200.15
$ws.Range("I132").Value = 952
$ws.Range("J132").Value = 1262.1875
$ws.Range("K132").Value = 8568
$ws.Range("L132").Value = 11359.6875
$ws.Range("M132").Value = -6038
$ws.Range("N132").Value = -16419.6875

$ws.Range("H135").Value = 1065.0227
$ws.Range("I135").Value = 250.22223
$ws.Range("K135").Value = 2252.00007
$ws.Range("M135").Value = 282.9999299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3863.1628
$ws.Range("I70").Value = 3912.8
$ws.Range("J70").Value = 3646
$ws.Range("K70").Value = 3912.8
$ws.Range("L70").Value = 3646
$ws.Range("M70").Value = -3642.8
$ws.Range("N70").Value = -4186

$ws.Range("H73").Value = 3863.1628
$ws.Range("I73").Value = 3912.8
$ws.Range("J73").Value = 3646
$ws.Range("K73").Value = 3912.8
$ws.Range("L73").Value = 3646
$ws.Range("M73").Value = -2976.8
$ws.Range("N73").Value = -5518

$ws.Range("H132").Value = 1660.1818
$ws.Range("I132").Value = 1649.5625
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4948.6875
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2418.6875
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3199.5881
$ws.Range("I40").Value = 2670.2307
$ws.Range("J40").Value = 4920
$ws.Range("K40").Value = 2670.2307
$ws.Range("L40").Value = 4920
$ws.Range("M40").Value = -2534.2307
$ws.Range("N40").Value = -5192

$ws.Range("H130").Value = 38660.555
$ws.Range("J130").Value = 38660.555
$ws.Range("L130").Value = 38660.555
$ws.Range("N130").Value = -48700.555

$ws.Range("H133").Value = 18607.875
$ws.Range("J133").Value = 18607.875
$ws.Range("L133").Value = 18607.875
$ws.Range("N133").Value = -23667.875

$ws.Range("H134").Value = 27639.125
$ws.Range("J134").Value = 33518.832
$ws.Range("L134").Value = 33518.832
$ws.Range("N134").Value = -43658.832

$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws.Range("H136").Value = 2154.253
$ws.Range("I136").Value = 1316.2142
$ws.Range("J136").Value = 3892.4075
$ws.Range("K136").Value = 3948.6426
$ws.Range("L136").Value = 11677.2225
$ws.Range("M136").Value = -1398.6426
$ws.Range("N136").Value = -16777.2225

$ws.Range("H140").Value = 48000
$ws.Range("J140").Value = 48000
$ws.Range("L140").Value = 48000
$ws.Range("N140").Value = -58360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 61285.8
$ws.Range("J46").Value = 51607.25
$ws.Range("L46").Value = 51607.25
$ws.Range("N46").Value = -52069.25

$ws.Range("H107").Value = 246.55556
$ws.Range("I107").Value = 166.625
$ws.Range("J107").Value = 362.81818
$ws.Range("K107").Value = 499.875
$ws.Range("L107").Value = 1088.45454
$ws.Range("M107").Value = 1420.125
$ws.Range("N107").Value = -4928.45454

$ws.Range("H123").Value = 25843.75
$ws.Range("J123").Value = 25843.75
$ws.Range("L123").Value = 25843.75
$ws.Range("N123").Value = -35643.75

$ws.Range("H132").Value = 750.26086
$ws.Range("I132").Value = 594.65717
$ws.Range("J132").Value = 1245.3636
$ws.Range("K132").Value = 1783.97151
$ws.Range("L132").Value = 3736.0908
$ws.Range("M132").Value = 746.0284900000001
$ws.Range("N132").Value = -8796.0908

$ws.Range("H134").Value = 61285.8
$ws.Range("J134").Value = 51607.25
$ws.Range("L134").Value = 154821.75
$ws.Range("N134").Value = -159891.75

$ws.Range("H136").Value = 416.72726
$ws.Range("I136").Value = 322.18182
$ws.Range("J136").Value = 700.36365
$ws.Range("K136").Value = 966.54546
$ws.Range("L136").Value = 2101.09095
$ws.Range("M136").Value = 1583.45454
$ws.Range("N136").Value = -7201.09095
